$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# Update the man.dm (column C) values for the remaining "digested biomass" rows
# from 5.1% to 5.9% dry matter
$ws.Range("C2:C5").Value = 5.9

# Remove the now-obsolete duplicate rows (previously 6.9% DM variants, rows 6-9)
$ws.Range("A6:D9").EntireRow.Delete()

# Match the cursor position left behind by the author's edit session
$ws.Range("F14").Select() | Out-Null
